function Get-ParaForRange($doc, $rng) {
    $all = $doc.Paragraphs
    for ($i = 1; $i -le $all.Count; $i++) {
        $pp = $all.Item($i)
        if ($pp.Range.Start -le $rng.Start -and $pp.Range.End -ge $rng.End) {
            return $pp
        }
    }
    return $null
}

$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$lsq = [char]0x2018
$rsq = [char]0x2019

# ---------------------------------------------------------------------------
# 1) "When you create a sponsorship..." paragraph (use case 005 Description):
#    add rPr (theme fonts + italic off) to paragraph mark and to each of the
#    3 existing runs; text stays the same.
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("When you create a sponsorship*it.", $false, $false, $true, $false, $false, $true, 1, $false, "", 0)
$p = Get-ParaForRange $d $rng
$full = $p.Range

$xml = $pkgOpen + '<w:p><w:pPr><w:pStyle w:val="Notes"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:i w:val="0"/></w:rPr></w:pPr>' `
    + '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:i w:val="0"/></w:rPr><w:t xml:space="preserve">When you create a sponsorship providing an expired credit card, the system must show the following message: ' + $lsq + 'Credit card must not have expired' + $rsq + ', but the system </w:t></w:r>' `
    + '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:i w:val="0"/></w:rPr><w:t>creates</w:t></w:r>' `
    + '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:i w:val="0"/></w:rPr><w:t xml:space="preserve"> it.</w:t></w:r>' `
    + '</w:p>' + $pkgClose
$full.InsertXML($xml)

# ---------------------------------------------------------------------------
# 2) Insert new paragraph after the "Results" that follows use case 005,
#    containing "The testers didn't report this bug" + the _GoBack bookmark
#    (moved here from the use case 006 description paragraph), underlined
#    paragraph mark.
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("Results", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p = Get-ParaForRange $d $rng
$insertPoint = $d.Range($p.Range.End - 1, $p.Range.End - 1)

$xml = $pkgOpen + '<w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' `
    + '<w:r><w:t>The testers didn' + $rsq + 't report this bug</w:t></w:r>' `
    + '<w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/>' `
    + '</w:p>' + $pkgClose
$insertPoint.InsertXML($xml)

# ---------------------------------------------------------------------------
# 3) Merge the 3 runs of "Bug in" + " use case 006: Cancel" + " a sponsorship"
#    into a single run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("Bug in use case 006*sponsorship", $false, $false, $true, $false, $false, $true, 1, $false, "", 0)
$p = Get-ParaForRange $d $rng
$full = $p.Range

$xml = $pkgOpen + '<w:p><w:pPr><w:pStyle w:val="Ttulo1"/></w:pPr>' `
    + '<w:r><w:t>Bug in use case 006: Cancel a sponsorship</w:t></w:r>' `
    + '</w:p>' + $pkgClose
$full.InsertXML($xml)

# ---------------------------------------------------------------------------
# 4) "When you cancel a sponsorship..." paragraph (use case 006 Description):
#    merge all runs into one, drop the _GoBack bookmark (moved in step 2),
#    add rPr (theme fonts + italic off).
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("When you cancel*cancel it.", $false, $false, $true, $false, $false, $true, 1, $false, "", 0)
$p = Get-ParaForRange $d $rng
$full = $p.Range

$xml = $pkgOpen + '<w:p><w:pPr><w:pStyle w:val="Notes"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:i w:val="0"/></w:rPr></w:pPr>' `
    + '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:i w:val="0"/></w:rPr><w:t>When you cancel a sponsorship providing an expired credit card, the system must return the list of sponsorship created with the sponsorship cancelled, but the system doesn' + $rsq + 't cancel it.</w:t></w:r>' `
    + '</w:p>' + $pkgClose
$full.InsertXML($xml)

# ---------------------------------------------------------------------------
# 5) The empty paragraph right after the use case 006 "Results" gets text
#    "The testers reported this bug."
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("Results", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$null = $rng.Find.Execute("Results", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p = Get-ParaForRange $d $rng
$emptyP = $p.Next()
$full = $emptyP.Range

$xml = $pkgOpen + '<w:p><w:r><w:t>The testers reported this bug.</w:t></w:r></w:p>' + $pkgClose
$full.InsertXML($xml)

# ---------------------------------------------------------------------------
# 6) "When you write or edit a review..." paragraph (use case 009
#    Description): split off the trailing period into its own run (without
#    rPr) and add rPr (theme fonts + italic off) to the first run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("When you write or edit a review*flag it.", $false, $false, $true, $false, $false, $true, 1, $false, "", 0)
$p = Get-ParaForRange $d $rng
$full = $p.Range

$xml = $pkgOpen + '<w:p><w:pPr><w:pStyle w:val="Notes"/></w:pPr>' `
    + '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:i w:val="0"/></w:rPr><w:t>When you write or edit a review with contains spam words, the system doesn' + $rsq + 't flag it</w:t></w:r>' `
    + '<w:r><w:t>.</w:t></w:r>' `
    + '</w:p>' + $pkgClose
$full.InsertXML($xml)

# ---------------------------------------------------------------------------
# 7) Insert new paragraph after the final "Results" with the tester note,
#    underlined paragraph mark.
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("Results", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$null = $rng.Find.Execute("Results", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$null = $rng.Find.Execute("Results", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p = Get-ParaForRange $d $rng
$insertPoint = $d.Range($p.Range.End - 1, $p.Range.End - 1)

$xml = $pkgOpen + '<w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' `
    + '<w:r><w:t>The tester reported that the system doesn' + $rsq + 't show any sponsorship but they didn' + $rsq + 't detect the bug</w:t></w:r>' `
    + '</w:p>' + $pkgClose
$insertPoint.InsertXML($xml)
